$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet (tab) name to reflect new "through" date
$ws.Name = "Through 2021-10-06"

# Update the label cell for the October row
$ws.Range("A11").Value = "October (through 10-06)"

# Update October row (row 11) values
$ws.Range("B11").Value = 6
$ws.Range("C11").Value = 9
$ws.Range("D11").Value = 11
$ws.Range("E11").Value = 15
$ws.Range("F11").Value = 7
$ws.Range("G11").Value = 33
$ws.Range("H11").Value = 46

# Update Total row (row 12) values
$ws.Range("B12").Value = 232
$ws.Range("C12").Value = 438
$ws.Range("D12").Value = 638
$ws.Range("E12").Value = 563
$ws.Range("F12").Value = 429
$ws.Range("G12").Value = 934
$ws.Range("H12").Value = 1295
